$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 101-102; everything from the old row 101
# downward (through 232) shifts down to 103-234.
$ws.Rows("101:102").Insert()

# Populate the new row 101 ("Primera") with the new weekly price point.
$ws.Cells.Item(101, 1).Value = 1
$ws.Cells.Item(101, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(101, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(101, 4).Value = 44579
$ws.Cells.Item(101, 5).Value = 15
$ws.Cells.Item(101, 6).Value = 100114014
$ws.Cells.Item(101, 7).Value = "Betarraga"
$ws.Cells.Item(101, 8).Value = "Sin especificar"
$ws.Cells.Item(101, 9).Value = "Primera"
$ws.Cells.Item(101, 10).Value = 1000
$ws.Cells.Item(101, 11).Value = 350
$ws.Cells.Item(101, 12).Value = 400
$ws.Cells.Item(101, 13).Value = 375
$ws.Cells.Item(101, 14).Value = '$/paquete 4 unidades'
$ws.Cells.Item(101, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(101, 16).Value = 94
$ws.Cells.Item(101, 17).Value = 4
$ws.Cells.Item(101, 18).Value = "Hortaliza"

# Populate the new row 102 ("Segunda") with the new weekly price point.
$ws.Cells.Item(102, 1).Value = 1
$ws.Cells.Item(102, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(102, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(102, 4).Value = 44579
$ws.Cells.Item(102, 5).Value = 15
$ws.Cells.Item(102, 6).Value = 100114014
$ws.Cells.Item(102, 7).Value = "Betarraga"
$ws.Cells.Item(102, 8).Value = "Sin especificar"
$ws.Cells.Item(102, 9).Value = "Segunda"
$ws.Cells.Item(102, 10).Value = 1200
$ws.Cells.Item(102, 11).Value = 350
$ws.Cells.Item(102, 12).Value = 400
$ws.Cells.Item(102, 13).Value = 375
$ws.Cells.Item(102, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(102, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(102, 16).Value = 75
$ws.Cells.Item(102, 17).Value = 5
$ws.Cells.Item(102, 18).Value = "Hortaliza"
